$wb = $excel.ActiveWorkbook

# Add a new worksheet named "Sheet2" after the existing Sheet1
$sheet1 = $wb.Worksheets.Item("Sheet1")
$newSheet = $wb.Worksheets.Add($null, $sheet1)
$newSheet.Name = "Sheet2"

# Populate Sheet2 with data
$newSheet.Range("A1").Value = "usernamepa"
$newSheet.Range("B1").Value = "password"
$newSheet.Range("A2").Value = "admin"
$newSheet.Range("B2").Value = "admin"

# Select B2 on the new sheet and make it the active sheet
$newSheet.Activate()
$newSheet.Range("B2").Select()

# Update selection on Sheet1 to B2 as well
$sheet1.Range("B2").Select()

# Re-activate Sheet2 so it's the active tab in the saved workbook
$newSheet.Activate()
$newSheet.Range("B2").Select()
